$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.300.13"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "3.310.34"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'558.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'142.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.96%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.311.13"
$ws.Range("E8").Value = "  -2.90%  "
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "3.882.20"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("E15").Value = "  -4.75%  "
$ws.Range("D16").Value = "3.306.13"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "60.314.36"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").Value = "'6.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("D20").Value = "'14.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "'8.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").Value = "'375.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("E24").Value = "  -3.60%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "3.442.06"
$ws.Range("E26").Value = "  -3.86%  "
$ws.Range("E27").Value = "  -7.27%  "
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'7.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -2.75%  "
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("D34").Value = "'22.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("E35").Value = "  -4.04%  "
$ws.Range("E36").Value = "  -4.85%  "
$ws.Range("D37").Value = "'166.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("E39").Value = "  -5.53%  "
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").Value = "3.341.31"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'26.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -12.85%  "
$ws.Range("E42").Value = "  -5.18%  "
$ws.Range("D43").Value = "'42.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("E45").Value = "  -3.64%  "
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("D48").Value = "2.372.17"
$ws.Range("E48").Value = "  -6.74%  "
$ws.Range("D49").Value = "'0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  -5.38%  "
$ws.Range("D51").Value = "'21.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.12%  "
